# Rename the "_old"/"_new" suffixed column headers (row 1) to the new
# format-version based suffixes "_FV2310" / "_FV2404", turn the header
# row + data range into a real Excel Table (so structured references /
# filtering work), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename headers -----------------------------------------------
# Columns A..J used the "_old" suffix, columns L..U used the "_new"
# suffix (column K just holds the literal "diff" and stays untouched).
$headersFV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    # columns L..U == 12..21, column K (11) is "diff" and is left as-is
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# --- 2. Turn A1:U72 into an Excel Table -------------------------------
$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
